$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 105.0175368920161
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 105.94536616336812

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 106.54467735532351
$ws.Range("D3").Value = 106.18326688907649
$ws.Range("E3").Value = 105.70986751021746

# Update the selection to match the new reduced highlighted range
$ws.Range("B1:E3").Select()
